$wb = $excel.ActiveWorkbook

# --- Rename header labels on the existing two sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "PO Forecast"

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Data rows
$ws3.Cells.Item(2, 1).Value = 45508.99999999999
$ws3.Cells.Item(2, 2).Value = 2
$ws3.Cells.Item(2, 3).Value = 1.086308252416892
$ws3.Cells.Item(2, 4).Value = 2.580337530715027

$ws3.Cells.Item(3, 1).Value = 45522.99999999999
$ws3.Cells.Item(3, 2).Value = 2
$ws3.Cells.Item(3, 3).Value = 1.178496148279797
$ws3.Cells.Item(3, 4).Value = 2.755494609168774

$ws3.Cells.Item(4, 1).Value = 45613.99999999999
$ws3.Cells.Item(4, 2).Value = 3
$ws3.Cells.Item(4, 3).Value = 2.217639956501871
$ws3.Cells.Item(4, 4).Value = 3.728339494426986

$ws3.Cells.Item(5, 1).Value = 45641.99999999999
$ws3.Cells.Item(5, 2).Value = 3
$ws3.Cells.Item(5, 3).Value = 2.514221434293423
$ws3.Cells.Item(5, 4).Value = 4.084809374335504

$ws3.Cells.Item(6, 1).Value = 45648.99999999999
$ws3.Cells.Item(6, 2).Value = 3
$ws3.Cells.Item(6, 3).Value = 2.611737651079856
$ws3.Cells.Item(6, 4).Value = 4.117195606210609

$ws3.Cells.Item(7, 1).Value = 45655.99999999999
$ws3.Cells.Item(7, 2).Value = 3
$ws3.Cells.Item(7, 3).Value = 2.622178791492126
$ws3.Cells.Item(7, 4).Value = 4.121497194601235

$ws3.Cells.Item(8, 1).Value = 45662.99999999999
$ws3.Cells.Item(8, 2).Value = 3
$ws3.Cells.Item(8, 3).Value = 2.702644671912489
$ws3.Cells.Item(8, 4).Value = 4.293534146305087

$ws3.Cells.Item(9, 1).Value = 45669.99999999999
$ws3.Cells.Item(9, 2).Value = 4
$ws3.Cells.Item(9, 3).Value = 2.795015468773334
$ws3.Cells.Item(9, 4).Value = 4.264731301814704

$ws3.Cells.Item(10, 1).Value = 45676.99999999999
$ws3.Cells.Item(10, 2).Value = 4
$ws3.Cells.Item(10, 3).Value = 2.848164266294901
$ws3.Cells.Item(10, 4).Value = 4.411500274766896

$ws3.Cells.Item(11, 1).Value = 45683.99999999999
$ws3.Cells.Item(11, 2).Value = 4
$ws3.Cells.Item(11, 3).Value = 2.961935806608567
$ws3.Cells.Item(11, 4).Value = 4.478063271184781

$ws3.Cells.Item(12, 1).Value = 45690.99999999999
$ws3.Cells.Item(12, 2).Value = 4
$ws3.Cells.Item(12, 3).Value = 2.986843722165821
$ws3.Cells.Item(12, 4).Value = 4.558710824934223

$ws3.Cells.Item(13, 1).Value = 45697.99999999999
$ws3.Cells.Item(13, 2).Value = 4
$ws3.Cells.Item(13, 3).Value = 3.077273350616812
$ws3.Cells.Item(13, 4).Value = 4.629268977136173

# --- Copy formatting from the source sheets so the new sheet matches ---
# Header style (bold, centered, bordered) -> row 1, columns A:D
$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Date style (custom date/time number format) -> column A, rows 2:13
$ws1.Range("A2").Copy()
$ws3.Range("A2:A13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Re-select A1 on the new sheet to match a fresh sheet view
$ws3.Range("A1").Select()
